$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows for accounts ASSAKO (004450724), BRUNO (004452912)
# and ERICA (004260002). Find rows by account number in column A and
# delete them (highest row index first so earlier deletions don't shift
# the row numbers of rows not yet processed).
$accountsToDelete = @("004450724", "004452912", "004260002")
$rowsToDelete = @()
foreach ($acct in $accountsToDelete) {
    $found = $ws.Columns(1).Find($acct)
    if ($found) {
        $rowsToDelete += $found.Row
    }
}
$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# Update JULIA's (004265173) balance from 30000 to 1000
$juliaCell = $ws.Columns(1).Find("004265173")
if ($juliaCell) {
    $ws.Cells.Item($juliaCell.Row, 3).Value = 1000
}

# Re-sort the data (excluding header row 1 and the trailing blank /
# filter-notes rows) in descending order by the Saldo column (C).
$lastRow = $ws.UsedRange.Rows.Count
$dataLastRow = $lastRow - 2
$sortRange = $ws.Range("A2:C" + $dataLastRow)
$keyRange = $ws.Range("C2:C" + $dataLastRow)
$sortRange.Sort($keyRange, 2)
